$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (2022-01-05): it used to be the most-recent row, so its
# numeric-looking figures were kept as text. Now that a newer row (10)
# exists, these settle into real numbers - matching the pattern already
# used by the older rows 6-8.
$ws.Range("B9").Value = 56348
$ws.Range("C9").Value = -2706
$ws.Range("E9").Value = 0

# --- Row 10: new data row for 2022-01-06 ("GoodInfo_v2 - 2022-01-06
# unfinished"), appended using the same text-formatted pattern row 9 used
# while it was still the newest entry.
$ws.Range("A10").Formula = "'2022-01-06"
$ws.Range("B10").Formula = "'56348.0"
$ws.Range("C10").Formula = "'-2074.0"
$ws.Range("D10").Formula = "'-3.68%"
$ws.Range("E10").Formula = "'0"

# Trailing (currently blank) performance columns, present on every row.
$ws.Range("F10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
